# Re-process the metadata for the "municipio-nombre" (column D) and
# "tipo-de-edificio" (column I) dimensions with the newly curated
# dimension/measure classification.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D = municipio-nombre metadata block (rows 2-4)
# It becomes a dimension (sdmx-dimension:refArea / dim) with a URI-Municipio
# mapping column, instead of a measure with xsd:int type.
$ws.Range("D2").Value = "sdmx-dimension:refArea"
$ws.Range("D3").Value = "dim"
$ws.Range("D4").Value = "URI-Municipio"

# Column I = tipo-de-edificio metadata block (rows 2-5)
# It becomes a measure (iaest-measure:tipo-de-edificio / medida / xsd:int)
# instead of a dimension, and no longer needs an external mapping file.
$ws.Range("I2").Value = "iaest-measure:tipo-de-edificio"
$ws.Range("I3").Value = "medida"
$ws.Range("I4").Value = "xsd:int"
$ws.Range("I5").Clear()
